$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45186
}

$hyperCols = @(19, 20, 21, 22, 23, 24, 25)
for ($r = 2; $r -le $lastRow; $r++) {
    $caseId = $ws.Cells.Item($r, 1).Value2
    foreach ($c in $hyperCols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -and $f.StartsWith("=HYPERLINK(")) {
            $inner = $f.Substring(11, $f.Length - 12)
            if (-not $inner.Contains(",")) {
                $cell.Formula = "=HYPERLINK(" + $inner + ", """ + $caseId + """)"
            }
        }
    }
}
